$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 210, shifting existing rows 210..258 down to 211..259
$ws.Rows("210:210").Insert()

# Populate the newly inserted row 210 with the new record
$ws.Range("A210").Value = 3
$ws.Range("B210").Value = "Femacal de La Calera"
$ws.Range("C210").Value = "Coquimbo"
$ws.Range("D210").Value = 44642
$ws.Range("E210").Value = 5
$ws.Range("F210").Value = 100112001
$ws.Range("G210").Value = "Berenjena"
$ws.Range("H210").Value = "Sin especificar"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 125
$ws.Range("K210").Value = 9000
$ws.Range("L210").Value = 9500
$ws.Range("M210").Value = 9260
$ws.Range("N210").Value = "`$/caja 60 unidades"
$ws.Range("O210").Value = "Región Metropolitana"
$ws.Range("P210").Value = 154
$ws.Range("Q210").Value = 60
$ws.Range("R210").Value = "Hortaliza"
